# The "Generative process" notes in the PRODUCT CODE column (D13:D15) are
# cleared out (replaced with a single blank space, matching the style
# already used elsewhere in the sheet), and the previously-empty
# PRICE PER UNIT (CHF) cells for those same rows (F13:F15) are likewise set
# to a single blank space. Once "Generative process" is no longer referenced
# anywhere, it drops out of the shared-strings table and everything that
# sorted after it shifts down by one index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

$ws.Range("D13").Value = " "
$ws.Range("F13").Value = " "

$ws.Range("D14").Value = " "
$ws.Range("F14").Value = " "

$ws.Range("D15").Value = " "
$ws.Range("F15").Value = " "

# Move/record the active selection onto F14, as in the final saved file.
$excel.Goto($ws.Range("F14"))
